$wb = $excel.ActiveWorkbook

# The new "Croatia" sheet mirrors the existing "Slovakia" sheet's layout
# (same columns/styles/merged cells), so copy it to the end of the
# workbook (after "Spain") and then adjust its content.
$slovakia = $wb.Worksheets.Item("Slovakia")
$spain = $wb.Worksheets.Item("Spain")
$slovakia.Copy($null, $spain)

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Market name + ticket reference for the new market.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2473"

# Restore Slovakia's own selection to "select all" (whole sheet), as it
# was left after the copy/paste operation performed to create Croatia.
$slovakia.Cells.Select()

# Make the new Croatia sheet the active tab/selection, matching the
# cursor being left on B4 after filling it in.
$croatia.Activate()
$croatia.Range("B4").Select()
